$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the "Unnamed: 0" header text from A1 (index column had no name).
#    Assigning a lone quote-prefix keeps the cell a (now empty) text cell
#    rather than turning it into a numeric blank.
$ws.Range("A1").Value = "'"

# 2. Remove the bold/border/centered formatting that was applied to the header
#    row (A1:AO1), reverting those cells to the default/unstyled format.
$ws.Range("A1:AO1").ClearFormats()

# 3. Corrected data values across the pre/post/total fixation metric rows
#    (rows 3-7) to account for a data-cleaning fix.
$ws.Range("F3").Value = 8
$ws.Range("H3").Value = 52
$ws.Range("I3").Value = 48
$ws.Range("J3").Value = 18
$ws.Range("L3").Value = 16
$ws.Range("O3").Value = 58
$ws.Range("P3").Value = 28
$ws.Range("Q3").Value = 48
$ws.Range("T3").Value = 24
$ws.Range("Y3").Value = 4
$ws.Range("Z3").Value = 7
$ws.Range("AA3").Value = 7
$ws.Range("AB3").Value = 52
$ws.Range("AG3").Value = 5
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 5
$ws.Range("AO3").Value = 2
$ws.Range("F4").Value = 10
$ws.Range("H4").Value = 338
$ws.Range("I4").Value = 205
$ws.Range("J4").Value = 29
$ws.Range("L4").Value = 24
$ws.Range("O4").Value = 261
$ws.Range("P4").Value = 55
$ws.Range("Q4").Value = 153
$ws.Range("T4").Value = 45
$ws.Range("Y4").Value = 5
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 248
$ws.Range("AG4").Value = 8
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 10
$ws.Range("AO4").Value = 4
$ws.Range("F5").Value = 2084.6
$ws.Range("H5").Value = 87328.31
$ws.Range("I5").Value = 52925.27
$ws.Range("J5").Value = 5604.74
$ws.Range("L5").Value = 5271.22
$ws.Range("O5").Value = 66956.61
$ws.Range("P5").Value = 13496.64
$ws.Range("Q5").Value = 38374.95
$ws.Range("T5").Value = 9825.8
$ws.Range("Y5").Value = 1084.39
$ws.Range("Z5").Value = 1969.02
$ws.Range("AA5").Value = 7407.87
$ws.Range("AB5").Value = 75845.34
$ws.Range("AG5").Value = 3286.99
$ws.Range("AJ5").Value = 4237.1
$ws.Range("AK5").Value = 2569.26
$ws.Range("AO5").Value = 734.18
$ws.Range("B6").Value = 0.05
$ws.Range("C6").Value = 0.17
$ws.Range("E6").Value = 0.05
$ws.Range("F6").Value = 0.65
$ws.Range("G6").Value = 6.14
$ws.Range("H6").Value = 27.32
$ws.Range("I6").Value = 16.56
$ws.Range("J6").Value = 1.75
$ws.Range("K6").Value = 3.78
$ws.Range("L6").Value = 1.65
$ws.Range("N6").Value = 1.74
$ws.Range("O6").Value = 20.95
$ws.Range("P6").Value = 4.22
$ws.Range("Q6").Value = 12.01
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 7.14
$ws.Range("T6").Value = 3.07
$ws.Range("U6").Value = 2.27
$ws.Range("V6").Value = 3.4
$ws.Range("W6").Value = 0.88
$ws.Range("X6").Value = 0.58
$ws.Range("Y6").Value = 0.34
$ws.Range("Z6").Value = 0.62
$ws.Range("AA6").Value = 2.32
$ws.Range("AB6").Value = 23.73
$ws.Range("AC6").Value = 0.14
$ws.Range("AD6").Value = 1.01
$ws.Range("AE6").Value = 0.48
$ws.Range("AF6").Value = 0.96
$ws.Range("AG6").Value = 1.03
$ws.Range("AH6").Value = 0.52
$ws.Range("AI6").Value = 0.17
$ws.Range("AJ6").Value = 1.33
$ws.Range("AK6").Value = 0.8
$ws.Range("AL6").Value = 0.49
$ws.Range("AM6").Value = 0.79
$ws.Range("AN6").Value = 0.45
$ws.Range("AO6").Value = 0.23
$ws.Range("F7").Value = 208.46
$ws.Range("H7").Value = 258.37
$ws.Range("I7").Value = 258.17
$ws.Range("J7").Value = 193.27
$ws.Range("L7").Value = 219.63
$ws.Range("O7").Value = 256.54
$ws.Range("P7").Value = 245.39
$ws.Range("Q7").Value = 250.82
$ws.Range("T7").Value = 218.35
$ws.Range("Y7").Value = 216.88
$ws.Range("Z7").Value = 179
$ws.Range("AA7").Value = 352.76
$ws.Range("AB7").Value = 305.83
$ws.Range("AG7").Value = 410.87
$ws.Range("AJ7").Value = 201.77
$ws.Range("AK7").Value = 256.93
$ws.Range("AO7").Value = 183.55

# 4. Drop the trailing fully-blank rows (10-14) that were left over from the
#    export; the used range should now end at row 9.
$ws.Rows("10:14").Delete()
